# Append new daily data rows (98-106) to the bottom of the swap table on Sheet1,
# extending the used range from A1:P97 to A1:P106.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(44340, 631, 2401, 1050, 6177, 967, 2188, 1000, 28500, 0, 0, 0, 3674, 356, 1073, 42939),
    @(44341, 415, 2301,  672, 5927, 218, 2153,    0, 28500, 0, 0, 0, 3674, 226, 1115, 42554),
    @(44342, 580, 2881,  450, 5177,  20, 2153, 1000, 28500, 0, 0, 0, 3674,  84, 1166, 42385),
    @(44343, 785, 2761,  750, 5027, 441, 2163,    0, 28500, 0, 0, 0, 3674, 463, 1153, 42125),
    @(44344, 440, 2851, 1144, 5121, 525, 2171, 1000, 28500, 0, 0, 0, 3674,  13, 1142, 42317),
    @(44347,   0, 2851,  286, 5407,   0, 2171,    0, 28500, 0, 0, 0, 3674,   0, 1142, 42603),
    @(44348, 300, 2736,  900, 5557, 214, 2167, 1500, 28500, 0, 0, 0, 3674, 123, 1039, 42634),
    @(44349, 755, 2911, 1627, 6157,  50, 2198, 1000, 28500, 0, 0, 0, 3674,  85, 1040, 43440),
    @(44350, 866, 2992, 1261, 6518, 439, 2196,    0, 28500, 0, 0, 0, 3674, 509, 1087, 43880)
)

$startRow = 98
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
    # Column A carries the sheet's date format (yyyy-mm-dd), same as every
    # preceding row in the table.
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd"
}
